$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("H29").Value = 399.5
$ws.Range("J29").Value = 400
$ws.Range("L29").Value = 1200
$ws.Range("N29").Value = -1762
$ws.Range("H34").Value = 2415.25
$ws.Range("I34").Value = 2415.25
$ws.Range("K34").Value = 2415.25
$ws.Range("M34").Value = -2212.25
$ws.Range("H36").Value = 2415.25
$ws.Range("I36").Value = 2415.25
$ws.Range("K36").Value = 2415.25
$ws.Range("M36").Value = -1700.25
$ws.Range("H38").Value = 109.1
$ws.Range("I38").Value = 121.111115
$ws.Range("J38").Value = 1
$ws.Range("K38").Value = 363.333345
$ws.Range("L38").Value = 3
$ws.Range("M38").Value = 8.666654999999992
$ws.Range("N38").Value = -747
$ws.Range("M21").ClearContents()
$ws.Range("M23").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3054.913
$ws.Range("I32").Value = 2953.024
$ws.Range("J32").Value = 4124.75
$ws.Range("K32").Value = 2953.024
$ws.Range("L32").Value = 4124.75
$ws.Range("M32").Value = -2666.024
$ws.Range("N32").Value = -4698.75
$ws.Range("H45").Value = 2166.5
$ws.Range("I45").Value = 1749.75
$ws.Range("K45").Value = 1749.75
$ws.Range("M45").Value = -1372.75
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("H97").Value = 1995.9
$ws.Range("I97").Value = 1326.6666
$ws.Range("J97").Value = 2999.75
$ws.Range("K97").Value = 1326.6666
$ws.Range("L97").Value = 2999.75
$ws.Range("M97").Value = -830.6666
$ws.Range("N97").Value = -3991.75
$ws.Range("L76").ClearContents()
$ws.Range("L79").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 19000.666
$ws.Range("J49").Value = 19000.666
$ws.Range("L49").Value = 19000.666
$ws.Range("N49").Value = -19478.666
$ws.Range("H99").Value = 1630
$ws.Range("I99").Value = 1645
$ws.Range("K99").Value = 1645
$ws.Range("M99").Value = -147
$ws.Range("H126").Value = 42780
$ws.Range("J126").Value = 42780
$ws.Range("L126").Value = 42780
$ws.Range("N126").Value = -52660

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2236.7144
$ws.Range("I58").Value = 1856.6
$ws.Range("K58").Value = 1856.6
$ws.Range("M58").Value = -1653.6
$ws.Range("H125").Value = 97619.25
$ws.Range("J125").Value = 97619.25
$ws.Range("L125").Value = 97619.25
$ws.Range("N125").Value = -102539.25
$ws.Range("H132").Value = 3346.375
$ws.Range("I132").Value = 3295.25
$ws.Range("K132").Value = 9885.75
$ws.Range("M132").Value = -7355.75
$ws.Range("H134").Value = 3297.4
$ws.Range("I134").Value = 2750
$ws.Range("K134").Value = 8250
$ws.Range("M134").Value = -5715
$ws.Range("H136").Value = 2236.7144
$ws.Range("I136").Value = 1856.6
$ws.Range("K136").Value = 5569.799999999999
$ws.Range("M136").Value = -3019.799999999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2327.5
$ws.Range("I39").Value = 346.57144
$ws.Range("J39").Value = 5100.8
$ws.Range("K39").Value = 1039.71432
$ws.Range("L39").Value = 15302.4
$ws.Range("M39").Value = -745.71432
$ws.Range("N39").Value = -15890.4
$ws.Range("H107").Value = 928.4545000000001
$ws.Range("I107").Value = 1022.5
$ws.Range("J107").Value = 815.6
$ws.Range("K107").Value = 3067.5
$ws.Range("L107").Value = 2446.8
$ws.Range("M107").Value = -1147.5
$ws.Range("N107").Value = -6286.8
$ws.Range("H132").Value = 1300
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 11882.909
$ws.Range("I24").Value = 356
$ws.Range("K24").Value = 356
$ws.Range("M24").Value = -183
$ws.Range("H132").Value = 3580.889
$ws.Range("I132").Value = 2760.611
$ws.Range("K132").Value = 8281.832999999999
$ws.Range("M132").Value = -5751.832999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 586.5
$ws.Range("I22").Value = 610.2222
$ws.Range("J22").Value = 556
$ws.Range("K22").Value = 610.2222
$ws.Range("L22").Value = 556
$ws.Range("M22").Value = -315.2222
$ws.Range("N22").Value = -1146
$ws.Range("H27").Value = 586.5
$ws.Range("I27").Value = 610.2222
$ws.Range("J27").Value = 556
$ws.Range("K27").Value = 610.2222
$ws.Range("L27").Value = 556
$ws.Range("M27").Value = -503.2222
$ws.Range("N27").Value = -770
$ws.Range("H46").Value = 3004.3635
$ws.Range("J46").Value = 4333.3335
$ws.Range("L46").Value = 4333.3335
$ws.Range("N46").Value = -4709.3335
$ws.Range("H55").Value = 620.1539
$ws.Range("I55").Value = 620.1539
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 620.1539
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = -447.1539
$ws.Range("H61").Value = 24657.791
$ws.Range("I61").Value = 20732.053
$ws.Range("K61").Value = 20732.053
$ws.Range("M61").Value = -20530.053
$ws.Range("H68").Value = 7293.5293
$ws.Range("I68").Value = 4785
$ws.Range("K68").Value = 4785
$ws.Range("M68").Value = -4036
$ws.Range("H71").Value = 7293.5293
$ws.Range("I71").Value = 4785
$ws.Range("K71").Value = 23925
$ws.Range("M71").Value = -20181
$ws.Range("H113").Value = 24657.791
$ws.Range("I113").Value = 20732.053
$ws.Range("K113").Value = 20732.053
$ws.Range("M113").Value = -18562.053
$ws.Range("H134").Value = 48333
$ws.Range("J134").Value = 48333
$ws.Range("L134").Value = 48333
$ws.Range("N134").Value = -58473
$ws.Range("H137").Value = 59000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 59000
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = 59000
$ws.Range("N137").Value = -69200
$ws.Range("M55").ClearContents()
$ws.Range("L137").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9055.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 9055.5
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = 9055.5
$ws.Range("N74").Value = -10927.5
$ws.Range("H77").Value = 9055.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 9055.5
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = 27166.5
$ws.Range("N77").Value = -36526.5
$ws.Range("H123").Value = 527499.5
$ws.Range("J123").Value = 527499.5
$ws.Range("L123").Value = 527499.5
$ws.Range("N123").Value = -537299.5
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360
$ws.Range("L74").ClearContents()
$ws.Range("L77").ClearContents()
